{"js": "// The document contains three short pieces of text that need updating:\n//   1. The student's name: \"BandiSreesaicharan\" -> \"G.AKSHAYA\"\n//   2. The batch number:   \"Batch-03\"           -> \"Batch-04\"\n//   3. The student id:     \"2403a54088\"         -> \"2403a54118\"\n// Everything else in the body text is unchanged (the underlying OOXML\n// merges/reflows some runs and drops now-stale proofing marks, but none of\n// that is observable/controllable through the Word JS API's text model).\n\nconst replacements = [\n  [\"BandiSreesaicharan\", \"G.AKSHAYA\"],\n  [\"Batch-03\", \"Batch-04\"],\n  [\"2403a54088\", \"2403a54118\"],\n];\n\nfor (const [find, replace] of replacements) {\n  const results = context.document.body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# The document contains three short pieces of text that need updating:\n#   1. The student's name: \"BandiSreesaicharan\" -> \"G.AKSHAYA\"\n#   2. The batch number:   \"Batch-03\"           -> \"Batch-04\"\n#   3. The student id:     \"2403a54088\"         -> \"2403a54118\"\n# Everything else in the body text is unchanged (the underlying OOXML\n# merges/reflows some runs and drops now-stale proofing marks, but none of\n# that is observable/controllable through the Word object model's text\n# model either).\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"BandiSreesaicharan\"; Replace = \"G.AKSHAYA\" },\n    @{ Find = \"Batch-03\";           Replace = \"Batch-04\" },\n    @{ Find = \"2403a54088\";         Replace = \"2403a54118\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Find\n    $find.Replacement.Text = $pair.Replace\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n}\n"}
